# Update "想去人数" (number of people interested) figures that changed
# between the previous and newly generated data snapshot.
# Affects both the "展览" sheet and the "全部类型" sheet (which mirrors it).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F8").Value = 1398
    $ws.Range("F11").Value = 1145
    $ws.Range("F15").Value = 37
}
